$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 640
$ws1.Range("F4").Value = 706
$ws1.Range("F5").Value = 584
$ws1.Range("F6").Value = 327
$ws1.Range("F7").Value = 2963
$ws1.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202407/R7iP9Iio1720170437964.jpeg"
$ws1.Range("F8").Value = 476
$ws1.Range("F9").Value = 8164
$ws1.Range("F12").Value = 55
$ws1.Range("F13").Value = 443
$ws1.Range("F14").Value = 53

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 640
$ws4.Range("F4").Value = 706
$ws4.Range("F5").Value = 584
$ws4.Range("F6").Value = 327
$ws4.Range("F9").Value = 2963
$ws4.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202407/R7iP9Iio1720170437964.jpeg"
$ws4.Range("F10").Value = 476
$ws4.Range("F12").Value = 8164
$ws4.Range("F15").Value = 55
$ws4.Range("F18").Value = 443
$ws4.Range("F19").Value = 53

$wb.Save()
